$d = $word.ActiveDocument

# Replace the old email domain with the new one.
$d.Content.Find.Execute("addison.boyer@umontana.edu", $true, $false, $false, $false, $false,
                         $true, 1, $false, "addison.boyer@mso.umt.edu", 2)
